$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated capital structure database values (applies identically to rows 2 and 3)
$values = @{
    "D"  = -0.098
    "G"  = -0.09304635761589404
    "H"  = -0.09304635761589404
    "I"  = -0.2059602649006622
    "J"  = -0.2059602649006622
    "K"  = -1.92
    "L"  = -0.06357615894039735
    "U"  = 1.27
    "V"  = 0.04686346863468634
    "W"  = -0.9948186528497409
    "X"  = 0.07951240380985054
    "Y"  = -1.074331056659591
    "Z"  = 1.080810249803164
    "AA" = -0.2226039653568105
    "AB" = 0.05953923903556256
    "AC" = -0.2821432043923731
    "AD" = 17.7
    "AE" = 0
    "AF" = 17.7
    "AG" = 16.43
    "AH" = 0.3950892857142857
    "AI" = 0.9404888416578108
    "AJ" = 0.3774408453939812
    "AK" = 0.9361823361823362
    "AL" = 1.12
    "AM" = 1.12
    "AN" = -3.094405594405595
    "AO" = -5.553571428571428
    "AP" = -2.872377622377623
    "AQ" = -5.553571428571428
}

foreach ($row in 2..3) {
    foreach ($col in $values.Keys) {
        $ws.Range($col + $row).Value = $values[$col]
    }
}

Write-Host "Done updating capital structure database."
